$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting from H1 (bold/border style) to new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in data values for column I (I0) and J (IF), rows 2-59
$iVals = @(9,9,8,8,9,9,8,7,8,8,6,7,7,8,10,8,8,7,7,8,8,8,6,7,8,9,5,6,9,10,9,8,8,8,8,9,9,8,9,9,9,8,5,10,1,2,1,1,1,5,6,7,5,6,4,1,1,1)
$jVals = @(9,9,9,9,9,9,9,7,8,8,7,7,7,8,10,9,8,8,7,8,9,8,6,7,8,9,5,6,9,10,9,8,9,8,8,9,9,8,9,10,9,8,7,10,1,3,3,4,5,8,7,7,5,9,6,3,3,2)

for ($r = 0; $r -lt $iVals.Length; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $iVals[$r]
    $ws.Cells.Item($row, 10).Value = $jVals[$r]
}
